# Se agregó tabla Localidades a BD
# Se agregó una tabla Localidades a la BD y se relaciono con la tabla Equipos.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Relate "Equipos" to the new "Localidades" table ---------------------
# The free-text "Localidad" attribute of Equipos becomes a foreign key
# ("LocalidadId", int) pointing at the new Localidades table.
$ws.Range("A40").Value = "LocalidadId"
$ws.Range("C40").Value = "FK"
$ws.Range("E40").Value = "int"

# --- 2. Add the new "Localidades" table right after "Contratos" -------------
# Rows 73-75 were blank; copy the formatting used by the other little
# "table" blocks on the sheet (header row style + field-row style, taken
# from the "Equipos" table at rows 37/38/40) onto the new rows, then fill
# in the table's header and its two fields (Id PK autoincrement int,
# Localidad nvarchar(40)).
$ws.Range("A37").Copy()
$ws.Range("A73:E73").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A38:E38").Copy()
$ws.Range("A74:E74").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A40:E40").Copy()
$ws.Range("A75:E75").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A73").Value = "Localidades"

$ws.Range("A74").Value = "Id"
$ws.Range("C74").Value = "PK"
$ws.Range("D74").Value = "Autoinc."
$ws.Range("E74").Value = "int"

$ws.Range("A75").Value = "Localidad"
$ws.Range("E75").Value = "nvarchar(40)"

$ws.Application.CutCopyMode = $false

# --- 3. Leave the sheet scrolled/selected on the new table -------------------
$ws.Range("B74").Select()
